$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in "Monthly EToF (unitless)s" header text (row 17, col A)
$ws.Range("A17").Value = "Monthly EToF (unitless)"

# Replace duplicated ET-rate values in the second table (rows 19-30) with
# the correct EToF (unitless) values
$ws.Range("C19").Value = 0.186
$ws.Range("D19").Value = 0.234
$ws.Range("E19").Value = 0.333
$ws.Range("F19").Value = 0.472
$ws.Range("G19").Value = 0.516
$ws.Range("H19").Value = 0.535
$ws.Range("I19").Value = 0.582
$ws.Range("J19").Value = 0.5659999999999999
$ws.Range("K19").Value = 0.527
$ws.Range("L19").Value = 0.5580000000000001
$ws.Range("M19").Value = 0.433
$ws.Range("N19").Value = 0.256
$ws.Range("C20").Value = 0.064
$ws.Range("D20").Value = 0.101
$ws.Range("E20").Value = 0.126
$ws.Range("F20").Value = 0.062
$ws.Range("G20").Value = 0.054
$ws.Range("H20").Value = 0.074
$ws.Range("I20").Value = 0.057
$ws.Range("J20").Value = 0.064
$ws.Range("K20").Value = 0.066
$ws.Range("L20").Value = 0.126
$ws.Range("M20").Value = 0.119
$ws.Range("N20").Value = 0.12
$ws.Range("C21").Value = 0.097
$ws.Range("D21").Value = 0.124
$ws.Range("E21").Value = 0.259
$ws.Range("F21").Value = 0.398
$ws.Range("G21").Value = 0.593
$ws.Range("H21").Value = 0.801
$ws.Range("I21").Value = 0.905
$ws.Range("J21").Value = 0.84
$ws.Range("K21").Value = 0.695
$ws.Range("L21").Value = 0.579
$ws.Range("M21").Value = 0.369
$ws.Range("N21").Value = 0.143
$ws.Range("C22").Value = 0.07199999999999999
$ws.Range("D22").Value = 0.07000000000000001
$ws.Range("E22").Value = 0.08400000000000001
$ws.Range("F22").Value = 0.064
$ws.Range("G22").Value = 0.08
$ws.Range("H22").Value = 0.096
$ws.Range("I22").Value = 0.077
$ws.Range("J22").Value = 0.078
$ws.Range("K22").Value = 0.095
$ws.Range("L22").Value = 0.106
$ws.Range("M22").Value = 0.112
$ws.Range("N22").Value = 0.09
$ws.Range("C23").Value = 0.184
$ws.Range("D23").Value = 0.223
$ws.Range("E23").Value = 0.312
$ws.Range("F23").Value = 0.426
$ws.Range("G23").Value = 0.551
$ws.Range("H23").Value = 0.636
$ws.Range("I23").Value = 0.657
$ws.Range("J23").Value = 0.63
$ws.Range("K23").Value = 0.627
$ws.Range("L23").Value = 0.569
$ws.Range("M23").Value = 0.445
$ws.Range("N23").Value = 0.244
$ws.Range("C24").Value = 0.102
$ws.Range("D24").Value = 0.107
$ws.Range("E24").Value = 0.109
$ws.Range("F24").Value = 0.07199999999999999
$ws.Range("G24").Value = 0.066
$ws.Range("H24").Value = 0.079
$ws.Range("I24").Value = 0.082
$ws.Range("J24").Value = 0.08599999999999999
$ws.Range("K24").Value = 0.091
$ws.Range("L24").Value = 0.103
$ws.Range("M24").Value = 0.096
$ws.Range("N24").Value = 0.07099999999999999
$ws.Range("C25").Value = 0.167
$ws.Range("D25").Value = 0.206
$ws.Range("E25").Value = 0.335
$ws.Range("F25").Value = 0.452
$ws.Range("G25").Value = 0.553
$ws.Range("H25").Value = 0.614
$ws.Range("I25").Value = 0.654
$ws.Range("J25").Value = 0.651
$ws.Range("K25").Value = 0.594
$ws.Range("L25").Value = 0.5629999999999999
$ws.Range("M25").Value = 0.449
$ws.Range("N25").Value = 0.297
$ws.Range("C26").Value = 0.068
$ws.Range("D26").Value = 0.076
$ws.Range("E26").Value = 0.098
$ws.Range("F26").Value = 0.066
$ws.Range("G26").Value = 0.04
$ws.Range("H26").Value = 0.06900000000000001
$ws.Range("I26").Value = 0.07099999999999999
$ws.Range("J26").Value = 0.077
$ws.Range("K26").Value = 0.07000000000000001
$ws.Range("L26").Value = 0.045
$ws.Range("M26").Value = 0.079
$ws.Range("N26").Value = 0.078
$ws.Range("C27").Value = 0.186
$ws.Range("D27").Value = 0.234
$ws.Range("E27").Value = 0.333
$ws.Range("F27").Value = 0.472
$ws.Range("G27").Value = 0.516
$ws.Range("H27").Value = 0.535
$ws.Range("I27").Value = 0.582
$ws.Range("J27").Value = 0.5659999999999999
$ws.Range("K27").Value = 0.527
$ws.Range("L27").Value = 0.5580000000000001
$ws.Range("M27").Value = 0.433
$ws.Range("N27").Value = 0.256
$ws.Range("C28").Value = 0.064
$ws.Range("D28").Value = 0.101
$ws.Range("E28").Value = 0.126
$ws.Range("F28").Value = 0.062
$ws.Range("G28").Value = 0.054
$ws.Range("H28").Value = 0.074
$ws.Range("I28").Value = 0.057
$ws.Range("J28").Value = 0.064
$ws.Range("K28").Value = 0.066
$ws.Range("L28").Value = 0.126
$ws.Range("M28").Value = 0.119
$ws.Range("N28").Value = 0.12
$ws.Range("C29").Value = 0.275
$ws.Range("D29").Value = 0.309
$ws.Range("E29").Value = 0.373
$ws.Range("F29").Value = 0.4
$ws.Range("G29").Value = 0.406
$ws.Range("H29").Value = 0.446
$ws.Range("I29").Value = 0.526
$ws.Range("J29").Value = 0.535
$ws.Range("K29").Value = 0.509
$ws.Range("L29").Value = 0.496
$ws.Range("M29").Value = 0.517
$ws.Range("N29").Value = 0.387
$ws.Range("C30").Value = 0.09
$ws.Range("D30").Value = 0.103
$ws.Range("E30").Value = 0.08699999999999999
$ws.Range("F30").Value = 0.036
$ws.Range("G30").Value = 0.036
$ws.Range("H30").Value = 0.052
$ws.Range("I30").Value = 0.051
$ws.Range("J30").Value = 0.055
$ws.Range("K30").Value = 0.054
$ws.Range("L30").Value = 0.049
$ws.Range("M30").Value = 0.06
$ws.Range("N30").Value = 0.07199999999999999
